# AfDD 2023 Annex Table Tab27 - data refresh + footnote encoding fix
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab27")

# Footnote A103: repair mojibake so PALOP/MERCOSUR names render with proper accents
$ws.Range("A103").Value = "Regional Economic Communities:CEN-SAD = ""Community of Sahel-Saharan States"";COMESA = ""Common Market for Eastern and Southern Africa"";EAC = ""East African Community"";ECCAS = ""Economic Community of Central African States"";ECOWAS = ""Economic Community of West African States"";IGAD = ""Intergovernmental Authority on Development"";SADC = ""Southern African Development Community"";UMA = ""Arab Maghreb Union"";PALOP = ""Países Africanos de Língua Oficial Portuguesa"";ASEAN = ""Association of Southeast Asian Nations"";MERCOSUR = ""Mercado Común del Sur"".EU27 = ""European Union (27 members)"".OECD = ""Organisation for Economic Co-operation and Development""."

# Recalculated figures (minor precision refresh)
$ws.Range("G67").Value = 758.719210156596
$ws.Range("N72").Value = 3.64174449612695
$ws.Range("H73").Value = 25686.237834625

# Row 97 - "Africa, Fragile States" - updated source data
$ws.Range("C97").Value = 93120.01662
$ws.Range("D97").Value = 8212.6333751374
$ws.Range("E97").Value = 405.235672683
$ws.Range("F97").Value = 521.829640886179
$ws.Range("G97").Value = 845.878247929779
$ws.Range("H97").Value = 7366.75512720763
$ws.Range("I97").Value = 49734.0115587128
$ws.Range("J97").Value = 281.575932898
$ws.Range("K97").Value = 571.474935779172
$ws.Range("L97").Value = 1282.09969335708
$ws.Range("M97").Value = 48451.9118653557
$ws.Range("N97").Value = 7.94646320746014
$ws.Range("O97").Value = 0.70394396331861
$ws.Range("P97").Value = 0.03403786547507
$ws.Range("Q97").Value = 0.04428301197498
$ws.Range("R97").Value = 0.07318143697628
$ws.Range("S97").Value = 0.63076252634233
$ws.Range("T97").Value = 4.14774356471227
$ws.Range("U97").Value = 0.02413896505068
$ws.Range("V97").Value = 0.05009437813465
$ws.Range("W97").Value = 0.11006418622681
$ws.Range("X97").Value = 4.03767937848546

# Row 98 - "ROW, Fragile States" - updated source data
$ws.Range("C98").Value = 154986.9198
$ws.Range("D98").Value = 14744.4263141411
$ws.Range("E98").Value = 1976.82581346481
$ws.Range("G98").Value = 6587.77330248503
$ws.Range("H98").Value = 8156.65301165608
$ws.Range("I98").Value = 83079.7189683261
$ws.Range("J98").Value = 5955.8849618
$ws.Range("L98").Value = 21107.889466406
$ws.Range("M98").Value = 61971.8295019201
$ws.Range("N98").Value = 12.686511270605
$ws.Range("O98").Value = 1.14579529558335
$ws.Range("P98").Value = 0.13180693375604
$ws.Range("R98").Value = 0.52359337443085
$ws.Range("S98").Value = 0.62220192115251
$ws.Range("T98").Value = 7.1301580404857
$ws.Range("U98").Value = 0.461884939688
$ws.Range("W98").Value = 1.95334595467436
$ws.Range("X98").Value = 5.17681208581133
